# Append the new survey response row (row 23) for "Robbie Sweetin" to Sheet1,
# matching the "Actualización SmartScore desde Streamlit (Robbie Sweetin)" export.
#
# All columns in this sheet (including the numeric-looking SmartScore columns,
# e.g. I/L/O/R/U/X/AA/AD/AG, and the blank Grupo_Experimental cell B23) are
# stored as TEXT in this dataset, except Edad (column D), which is a real
# number. Assigning a plain numeric-looking string straight to .Value would
# get auto-coerced to a Number by Excel, so every text cell is written with a
# leading apostrophe (forces literal text, just like typing '0.562 into a
# cell) and then the cell style is reset back to "Normal" so no quote-prefix
# formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddress, $text) {
    $ws.Range($cellAddress).Value = "'" + $text
    $ws.Range($cellAddress).Style = "Normal"
}

Set-TextValue "A23"  "Robbie Sweetin_20251202_131942"
Set-TextValue "B23"  ""
Set-TextValue "C23"  "Robbie Sweetin"
$ws.Range("D23").Value = 18
Set-TextValue "E23"  "Male"
Set-TextValue "F23"  "2025-12-02 13:19:42"
Set-TextValue "G23"  @'
{
  "portion": 0.6,
  "diet": 0.7142857142857143,
  "salt": 0.4,
  "fat": 0.8,
  "natural": 0.6,
  "convenience": 0.4,
  "price": 0.6
}
'@
Set-TextValue "H23"  "Nongshim Neoguri Spicy Seafood"
Set-TextValue "I23"  "0.562"
Set-TextValue "J23"  "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
Set-TextValue "K23"  "Maruchan Ramen Sabor Pollo"
Set-TextValue "L23"  "0.468"
Set-TextValue "M23"  "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"
Set-TextValue "N23"  "Nissin Chow Mein Teriyaki Beef"
Set-TextValue "O23"  "0.465"
Set-TextValue "P23"  "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"
Set-TextValue "Q23"  "Kraft Macaroni & Cheese Dinner"
Set-TextValue "R23"  "0.636"
Set-TextValue "S23"  "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
Set-TextValue "T23"  "Annie’s Shells & White Cheddar"
Set-TextValue "U23"  "0.581"
Set-TextValue "V23"  "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"
Set-TextValue "W23"  "Amy’s Macaroni & Cheese (frozen)"
Set-TextValue "X23"  "0.569"
Set-TextValue "Y23"  "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"
Set-TextValue "Z23"  "Wild Planet Wild Tuna Pasta Salad"
Set-TextValue "AA23" "0.706"
Set-TextValue "AB23" "Sabor fresco, buena proteína, saludable, porción algo pequeña"
Set-TextValue "AC23" "StarKist Chicken Creations (Chicken Salad)"
Set-TextValue "AD23" "0.582"
Set-TextValue "AE23" "Portátil, saludable, fácil, buena textura, sabor suave"
Set-TextValue "AF23" "Jack Link’s Beef Jerky Original"
Set-TextValue "AG23" "0.556"
Set-TextValue "AH23" "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"

# Multi-line cell content (G23) makes Excel auto-expand the row height; drop
# it back to an auto-fit (default) height so no explicit row height sticks.
$ws.Rows(23).AutoFit()
